{"js": "// Replace the date line and the 25 division problems with their updated\n// values, per the commit diff. Text values in this document are unique,\n// so a targeted search/replace for each old string is unambiguous.\nconst replacements = [\n  [\"2025-05-27 Tuesday\", \"2025-05-28 Wednesday\"],\n  [\"855\u00f72=\", \"305\u00f72=\"],\n  [\"753\u00f73=\", \"284\u00f72=\"],\n  [\"387\u00f72=\", \"807\u00f77=\"],\n  [\"179\u00f75=\", \"580\u00f78=\"],\n  [\"551\u00f76=\", \"710\u00f77=\"],\n  [\"212\u00f73=\", \"727\u00f75=\"],\n  [\"172\u00f77=\", \"224\u00f79=\"],\n  [\"156\u00f75=\", \"742\u00f74=\"],\n  [\"467\u00f73=\", \"951\u00f75=\"],\n  [\"731\u00f76=\", \"195\u00f72=\"],\n  [\"464\u00f78=\", \"607\u00f74=\"],\n  [\"470\u00f75=\", \"399\u00f79=\"],\n  [\"349\u00f78=\", \"620\u00f73=\"],\n  [\"570\u00f77=\", \"136\u00f74=\"],\n  [\"505\u00f77=\", \"189\u00f72=\"],\n  [\"675\u00f72=\", \"493\u00f79=\"],\n  [\"334\u00f72=\", \"966\u00f77=\"],\n  [\"978\u00f72=\", \"449\u00f78=\"],\n  [\"989\u00f76=\", \"316\u00f75=\"],\n  [\"426\u00f73=\", \"777\u00f73=\"],\n  [\"192\u00f72=\", \"332\u00f76=\"],\n  [\"506\u00f73=\", \"657\u00f72=\"],\n  [\"221\u00f79=\", \"832\u00f79=\"],\n  [\"712\u00f74=\", \"895\u00f75=\"],\n  [\"840\u00f76=\", \"984\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 division problems with their updated\n# values, per the commit diff. Text values in this document are unique,\n# so a targeted Find/Replace for each old string is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-05-27 Tuesday\", \"2025-05-28 Wednesday\"),\n    @(\"855\u00f72=\", \"305\u00f72=\"),\n    @(\"753\u00f73=\", \"284\u00f72=\"),\n    @(\"387\u00f72=\", \"807\u00f77=\"),\n    @(\"179\u00f75=\", \"580\u00f78=\"),\n    @(\"551\u00f76=\", \"710\u00f77=\"),\n    @(\"212\u00f73=\", \"727\u00f75=\"),\n    @(\"172\u00f77=\", \"224\u00f79=\"),\n    @(\"156\u00f75=\", \"742\u00f74=\"),\n    @(\"467\u00f73=\", \"951\u00f75=\"),\n    @(\"731\u00f76=\", \"195\u00f72=\"),\n    @(\"464\u00f78=\", \"607\u00f74=\"),\n    @(\"470\u00f75=\", \"399\u00f79=\"),\n    @(\"349\u00f78=\", \"620\u00f73=\"),\n    @(\"570\u00f77=\", \"136\u00f74=\"),\n    @(\"505\u00f77=\", \"189\u00f72=\"),\n    @(\"675\u00f72=\", \"493\u00f79=\"),\n    @(\"334\u00f72=\", \"966\u00f77=\"),\n    @(\"978\u00f72=\", \"449\u00f78=\"),\n    @(\"989\u00f76=\", \"316\u00f75=\"),\n    @(\"426\u00f73=\", \"777\u00f73=\"),\n    @(\"192\u00f72=\", \"332\u00f76=\"),\n    @(\"506\u00f73=\", \"657\u00f72=\"),\n    @(\"221\u00f79=\", \"832\u00f79=\"),\n    @(\"712\u00f74=\", \"895\u00f75=\"),\n    @(\"840\u00f76=\", \"984\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
